$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 32   Number  47"
$ws.Range("C9").Value = "Report Covering the Week  11/17/2025  Through  11/23/2025"

# --- Crime statistics table updates (rows 15-31) ---
$ws.Range("F15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 1
$ws.Range("F15").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("K23").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = 0
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 25
$ws.Range("J15").Value = 17
$ws.Range("K15").Value = 47.058823529411
$ws.Range("L15").Value = 150
$ws.Range("C16").Value = 1
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "***.*"
$ws.Range("A15").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 140
$ws.Range("I16").Value = 72
$ws.Range("K16").Value = -16.279069767441
$ws.Range("L16").Value = -1.369863013698
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 17
$ws.Range("H17").Value = 30.769230769230
$ws.Range("I17").Value = 219
$ws.Range("J17").Value = 185
$ws.Range("K17").Value = 18.378378378378
$ws.Range("L17").Value = -3.524229074889
$ws.Range("F15").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -37.5
$ws.Range("I18").Value = 67
$ws.Range("J18").Value = 92
$ws.Range("K18").Value = -27.173913043478
$ws.Range("L18").Value = 3.076923076923
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 200
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 34.782608695652
$ws.Range("I19").Value = 388
$ws.Range("J19").Value = 323
$ws.Range("K19").Value = 20.123839009287
$ws.Range("L19").Value = 0.779220779220
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -60
$ws.Range("J20").Value = 85
$ws.Range("K20").Value = -20
$ws.Range("L20").Value = -31.313131313131
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 30.769230769230
$ws.Range("F21").Value = 71
$ws.Range("G21").Value = 62
$ws.Range("H21").Value = 14.516129032258
$ws.Range("I21").Value = 839
$ws.Range("J21").Value = 791
$ws.Range("K21").Value = 6.068268015170
$ws.Range("L21").Value = -2.780996523754
$ws.Range("F23").Value = 1
$ws.Range("H23").Value = -66.666666666666
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -35.483870967741
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 121
$ws.Range("H24").Value = -20.661157024793
$ws.Range("I24").Value = 1209
$ws.Range("J24").Value = 1299
$ws.Range("K24").Value = -6.928406466512
$ws.Range("L24").Value = -9.505988023952
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -52.380952380952
$ws.Range("F25").Value = 69
$ws.Range("G25").Value = 86
$ws.Range("H25").Value = -19.767441860465
$ws.Range("I25").Value = 834
$ws.Range("J25").Value = 886
$ws.Range("K25").Value = -5.869074492099
$ws.Range("L25").Value = 6.513409961685
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 9
$ws.Range("E26").Value = 55.555555555555
$ws.Range("F26").Value = 35
$ws.Range("G26").Value = 34
$ws.Range("H26").Value = 2.941176470588
$ws.Range("I26").Value = 515
$ws.Range("J26").Value = 527
$ws.Range("K26").Value = -2.277039848197
$ws.Range("L26").Value = 10.991379310344
$ws.Range("F15").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 1
$ws.Range("F15").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("K23").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = 0
$ws.Range("I27").Value = 29
$ws.Range("J27").Value = 28
$ws.Range("K27").Value = 3.571428571428
$ws.Range("L27").Value = 38.095238095238
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 2
$ws.Range("F28").Value = 4
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -20
$ws.Range("I28").Value = 54
$ws.Range("J28").Value = 59
$ws.Range("K28").Value = -8.474576271186
$ws.Range("L28").Value = -6.896551724137
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("A15").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("A15").Copy()
$ws.Range("E31").PasteSpecial(-4122)
